$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update shared strings (route text) by replacing cell values wherever they occur.
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 17)  # Column Q
    $v = $cell.Value2
    if ($v -eq "purchase.route_warehouse0_buy") {
        $cell.Value = "purchase_stock.route_warehouse0_buy"
    } elseif ($v -eq "purchase.route_warehouse0_buy,stock.route_warehouse0_mto") {
        $cell.Value = "purchase_stock.route_warehouse0_buy,stock.route_warehouse0_mto"
    }
}

Write-Host "done"
